# Update the single data row and remove the extra rows (3-6),
# leaving only the header row and one data row (A1:B2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data row values
$ws.Range("A2").Value = "Hari"
$ws.Range("B2").Value = "31-10-1999"

# Delete the now-unused rows 3 through 6 (shift cells up / remove rows)
$ws.Range("A3:B6").EntireRow.Delete()

# Move the active selection to B2 to match the saved view state
$ws.Range("B2").Select()
